$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Copy formatting from the previous row (row 53) down onto row 54 first,
# so the new row matches existing styling (date format on column A, etc.)
$ws.Range("A53:E53").Copy()
$ws.Range("A54:E54").PasteSpecial(-4122)  # xlPasteFormats

# New row 54 data
$ws.Range("A54").Value = 45986
$ws.Range("B54").Value = 2025
$ws.Range("C54").Value = 2.46481303148316
$ws.Range("D54").Value = 2026
$ws.Range("E54").Value = 2.509429409292352
